$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3 each describe an observation record. The edit swaps the
# record data between row 2 and row 3 (same two records, rows exchanged),
# while the "Taxonsorteringsordning" (column B) gets fresh values that do
# not simply follow the swap.
#
# Use a temporary holding row (row 10, well below the used range) so the
# row contents -- including cell data types/number formats -- are moved
# with Copy instead of being retyped by hand (this keeps text-looking
# dates such as "2023-05-12" stored as text rather than being reinterpreted
# as date serial numbers).

$tempRow = 10

$ws.Range("A2:AY2").Copy($ws.Range("A" + $tempRow + ":AY" + $tempRow))
$ws.Range("A3:AY3").Copy($ws.Range("A2:AY2"))
$ws.Range("A" + $tempRow + ":AY" + $tempRow).Copy($ws.Range("A3:AY3"))
$ws.Range("A" + $tempRow + ":AY" + $tempRow).Clear()

# Fix up column B (Taxonsorteringsordning) with its new values.
$ws.Range("B2").Value = 56841
$ws.Range("B3").Value = 56894
